$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "GMHOR:0000002"
$ws.Range("C4:G4").ClearContents()
